$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 116, pushing the existing rows 116-126 down
# to 117-127 (dimension grows from A1:R126 to A1:R127).
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row 116 with the new weekly price entry.
$ws.Cells.Item(116, 1).Value = 11
$ws.Cells.Item(116, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(116, 3).Value = "Bíobío"
$ws.Cells.Item(116, 4).Value = 44491
$ws.Cells.Item(116, 5).Value = 8
$ws.Cells.Item(116, 6).Value = 100114001
$ws.Cells.Item(116, 7).Value = "Papa"
$ws.Cells.Item(116, 8).Value = "Asterix"
$ws.Cells.Item(116, 9).Value = "1a (guarda)"
$ws.Cells.Item(116, 10).Value = 2000
$ws.Cells.Item(116, 11).Value = 10500
$ws.Cells.Item(116, 12).Value = 11000
$ws.Cells.Item(116, 13).Value = 10750
$ws.Cells.Item(116, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(116, 15).Value = "Provincia de Arauco"
$ws.Cells.Item(116, 16).Value = 430
$ws.Cells.Item(116, 17).Value = 25
$ws.Cells.Item(116, 18).Value = "Hortaliza"
